# Update UHQ rates to reflect equations in paper (lines 58-63)
# Rewrites columns B:J for rows 3-22 (time steps 1-20) with recomputed
# susceptible/infected/UD/UR/HD/HR/QD/QR/died values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 20,9
$arr[0,0] = [double]"0.6610309099806798"
$arr[0,1] = [double]"0.1188067775532411"
$arr[0,2] = [double]"0.0001083189802580445"
$arr[0,3] = [double]"0.0107235790455464"
$arr[0,4] = [double]"4.061961759676667e-06"
$arr[0,5] = [double]"0.0004024156795688895"
$arr[0,6] = [double]"2.331262728070345e-05"
$arr[0,7] = [double]"0.002716118733894493"
$arr[0,8] = [double]"0.1180921373645059"
$arr[1,0] = [double]"0.6033558379651277"
$arr[1,1] = [double]"0.1429607227367725"
$arr[1,2] = [double]"0.0002478541818747085"
$arr[1,3] = [double]"0.01411739725565055"
$arr[1,4] = [double]"9.294531820301569e-06"
$arr[1,5] = [double]"0.0005370910581833254"
$arr[1,6] = [double]"5.295560199292023e-05"
$arr[1,7] = [double]"0.003078316154487901"
$arr[1,8] = [double]"0.1180958239391973"
$arr[2,0] = [double]"0.5359785339991505"
$arr[2,1] = [double]"0.1644335827547162"
$arr[2,2] = [double]"0.0004124671320663437"
$arr[2,3] = [double]"0.01699094641782112"
$arr[2,4] = [double]"1.546751745248789e-05"
$arr[2,5] = [double]"0.0006473095863210778"
$arr[2,6] = [double]"8.792782943499545e-05"
$arr[2,7] = [double]"0.003698622445585419"
$arr[2,8] = [double]"0.1181042489748803"
$arr[3,0] = [double]"0.4650191159006192"
$arr[3,1] = [double]"0.1841531339022194"
$arr[3,2] = [double]"0.0005982383516741455"
$arr[3,3] = [double]"0.01954671242480341"
$arr[3,4] = [double]"2.243393818778045e-05"
$arr[3,5] = [double]"0.0007451094171468894"
$arr[3,6] = [double]"0.0001273964138132682"
$arr[3,7] = [double]"0.004258755709809346"
$arr[3,8] = [double]"0.1181182641307136"
$arr[4,0] = [double]"0.3950488935544532"
$arr[4,1] = [double]"0.2025092421498664"
$arr[4,2] = [double]"0.0008024714385099327"
$arr[4,3] = [double]"0.02189480104361739"
$arr[4,4] = [double]"3.009267894412248e-05"
$arr[4,5] = [double]"0.000834857334227037"
$arr[4,6] = [double]"0.0001707883633715627"
$arr[4,7] = [double]"0.004772875656516981"
$arr[4,8] = [double]"0.1181385879582081"
$arr[5,0] = [double]"0.3292647662262045"
$arr[5,1] = [double]"0.2195514028011575"
$arr[5,2] = [double]"0.001023013328811494"
$arr[5,3] = [double]"0.02408144988626477"
$arr[5,4] = [double]"3.836299983043104e-05"
$arr[5,5] = [double]"0.0009183749060616591"
$arr[5,6] = [double]"0.0002176461459453748"
$arr[5,7] = [double]"0.005251236499678021"
$arr[5,8] = [double]"0.1181658474375374"
$arr[6,0] = [double]"0.2697048107417248"
$arr[6,1] = [double]"0.2351391741760236"
$arr[6,2] = [double]"0.001257830638966161"
$arr[6,3] = [double]"0.02611247634064593"
$arr[6,4] = [double]"4.716864896123105e-05"
$arr[6,5] = [double]"0.0009959179726567527"
$arr[6,6] = [double]"0.0002675376615731995"
$arr[6,7] = [double]"0.00569538614585585"
$arr[6,8] = [double]"0.1182005964295872"
$arr[7,0] = [double]"0.2174595847269764"
$arr[7,1] = [double]"0.2490511482660261"
$arr[7,2] = [double]"0.001504778390086337"
$arr[7,3] = [double]"0.02797113676370552"
$arr[7,4] = [double]"5.642918962823764e-05"
$arr[7,5] = [double]"0.001066865829623726"
$arr[7,6] = [double]"0.000320007096464036"
$arr[7,7] = [double]"0.006101824188918121"
$arr[7,8] = [double]"0.1182433197403573"
$arr[8,0] = [double]"0.1728730253175389"
$arr[8,1] = [double]"0.2610600687904289"
$arr[8,2] = [double]"0.001761505973083151"
$arr[8,3] = [double]"0.02963106546515542"
$arr[8,4] = [double]"6.605647399061816e-05"
$arr[8,5] = [double]"0.001130220062619847"
$arr[8,6] = [double]"0.0003745549405708796"
$arr[8,7] = [double]"0.006464859321873347"
$arr[8,8] = [double]"0.1182944293311354"
$arr[9,0] = [double]"0.1357337436702637"
$arr[9,1] = [double]"0.270980849234879"
$arr[9,2] = [double]"0.002025455888385338"
$arr[9,3] = [double]"0.03106519599611967"
$arr[9,4] = [double]"7.595459581445016e-05"
$arr[9,5] = [double]"0.001184948875676186"
$arr[9,6] = [double]"0.0004306377197764544"
$arr[9,7] = [double]"0.006778588842880892"
$arr[9,8] = [double]"0.1183542573351966"
$arr[10,0] = [double]"0.1054532548958411"
$arr[10,1] = [double]"0.2786970820657173"
$arr[10,2] = [double]"0.002293920142592403"
$arr[10,3] = [double]"0.03225146565545564"
$arr[10,4] = [double]"8.602200534721511e-05"
$arr[10,5] = [double]"0.001230207691529504"
$arr[10,6] = [double]"0.0004876799801776078"
$arr[10,7] = [double]"0.007038168474262355"
$arr[10,8] = [double]"0.1184230490317034"
$arr[11,0] = [double]"0.08122352199166959"
$arr[11,1] = [double]"0.284170907685431"
$arr[11,2] = [double]"0.002564126877078677"
$arr[11,3] = [double]"0.03317596642555665"
$arr[11,4] = [double]"9.615475789045036e-05"
$arr[11,5] = [double]"0.001265461529818443"
$arr[11,6] = [double]"0.0005450926966593927"
$arr[11,7] = [double]"0.007240518539350068"
$arr[11,8] = [double]"0.1185009577348645"
$arr[12,0] = [double]"0.06214598912384005"
$arr[12,1] = [double]"0.2874407401606637"
$arr[12,2] = [double]"0.002833336413742991"
$arr[12,3] = [double]"0.0338341199251522"
$arr[12,4] = [double]"0.0001062501155153621"
$arr[12,5] = [double]"0.001290531347417073"
$arr[12,6] = [double]"0.0006022936826127146"
$arr[12,7] = [double]"0.007384593654013793"
$arr[12,8] = [double]"0.1185880426250438"
$arr[13,0] = [double]"0.04732759066751309"
$arr[13,1] = [double]"0.2886110216189984"
$arr[13,2] = [double]"0.003098931882350984"
$arr[13,3] = [double]"0.03423041002072676"
$arr[13,4] = [double]"0.0001162099455881619"
$arr[13,5] = [double]"0.001305584805555648"
$arr[13,6] = [double]"0.0006587268482405278"
$arr[13,7] = [double]"0.007471333160687757"
$arr[13,8] = [double]"0.1186842698416305"
$arr[14,0] = [double]"0.03594391918243994"
$arr[14,1] = [double]"0.2878376812170164"
$arr[14,2] = [double]"0.003358495015587955"
$arr[14,3] = [double]"0.03437716828611342"
$arr[14,4] = [double]"0.0001259435630845483"
$arr[14,5] = [double]"0.001311090446540417"
$arr[14,6] = [double]"0.0007138783070026377"
$arr[14,7] = [double]"0.007503401035747628"
$arr[14,8] = [double]"0.1187895166536603"
$arr[15,0] = [double]"0.02727379174787814"
$arr[15,1] = [double]"0.2853122633939841"
$arr[15,2] = [double]"0.003609862372611912"
$arr[15,3] = [double]"0.03429284873353516"
$arr[15,4] = [double]"0.0001353698389729467"
$arr[15,5] = [double]"0.001307752024484692"
$arr[15,6] = [double]"0.000767288323116903"
$arr[15,7] = [double]"0.007484811232463836"
$arr[15,8] = [double]"0.1189035782044729"
$arr[16,0] = [double]"0.02071178317371192"
$arr[16,1] = [double]"0.2812468518289091"
$arr[16,2] = [double]"0.003851160942251559"
$arr[16,3] = [double]"0.03400014388624832"
$arr[16,4] = [double]"0.0001444185353344334"
$arr[16,5] = [double]"0.001296436522695202"
$arr[16,6] = [double]"0.0008185588772350957"
$arr[16,7] = [double]"0.007420516096245826"
$arr[16,8] = [double]"0.1190261761657504"
$arr[17,0] = [double]"0.01576579506419083"
$arr[17,1] = [double]"0.275861086836404"
$arr[17,2] = [double]"0.004080824658063207"
$arr[17,3] = [double]"0.03352419478069391"
$arr[17,4] = [double]"0.0001530309246773702"
$arr[17,5] = [double]"0.001278105593438077"
$arr[17,6] = [double]"0.0008673571758888849"
$arr[17,7] = [double]"0.007316013792821563"
$arr[17,8] = [double]"0.1191569686006572"
$arr[18,0] = [double]"0.01204590676195057"
$arr[18,1] = [double]"0.269371870112436"
$arr[18,2] = [double]"0.004297594871841122"
$arr[18,3] = [double]"0.03289104911371443"
$arr[18,4] = [double]"0.000161159807694042"
$arr[18,5] = [double]"0.001253756390368878"
$arr[18,6] = [double]"0.0009134157521072791"
$arr[18,7] = [double]"0.007177009138556126"
$arr[18,8] = [double]"0.1192955603878405"
$arr[19,0] = [double]"0.00924928192884648"
$arr[19,1] = [double]"0.2619858324171081"
$arr[19,2] = [double]"0.004500508459564494"
$arr[19,3] = [double]"0.03212643818012798"
$arr[19,4] = [double]"0.0001687690672336685"
$arr[19,5] = [double]"0.001224374555324059"
$arr[19,6] = [double]"0.0009565299379269465"
$arr[19,7] = [double]"0.007009143826920431"
$arr[19,8] = [double]"0.11944151366151"

$ws.Range("B3:J22").Value = $arr
